$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team-specific time-matrix probabilities (Greensboro_B)
# Values below replace the prior (generic) matrix with team-specific figures.

# Row 2
$ws.Range("B2").Value = 0.191044776119403
$ws.Range("C2").Value = 0.5970149253731343
$ws.Range("J2").Value = 0.01194029850746269
$ws.Range("P2").Value = 0.1313432835820895
$ws.Range("S2").Value = 0.06865671641791045

# Row 3
$ws.Range("B3").Value = 0.009389671361502348
$ws.Range("C3").Value = 0.04225352112676056
$ws.Range("J3").Value = 0.02347417840375587
$ws.Range("P3").Value = 0.7370892018779343
$ws.Range("S3").Value = 0.1877934272300469

# Row 4
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.2307692307692308

# Row 6
$ws.Range("B6").Value = 0.06772908366533864
$ws.Range("D6").Value = 0.00398406374501992
$ws.Range("F6").Value = 0.07171314741035857
$ws.Range("J6").Value = 0.302788844621514
$ws.Range("O6").Value = 0.04382470119521913
$ws.Range("Q6").Value = 0.1593625498007968
$ws.Range("R6").Value = 0.05976095617529881
$ws.Range("S6").Value = 0.2908366533864542

# Row 7
$ws.Range("B7").Value = 0.1046511627906977
$ws.Range("D7").Value = 0.005813953488372093
$ws.Range("F7").Value = 0.03488372093023256
$ws.Range("J7").Value = 0.186046511627907
$ws.Range("O7").Value = 0.01162790697674419
$ws.Range("Q7").Value = 0.1686046511627907
$ws.Range("R7").Value = 0.05813953488372093
$ws.Range("S7").Value = 0.4302325581395349

# Row 8
$ws.Range("B8").Value = 0.09606986899563319
$ws.Range("D8").Value = 0.01091703056768559
$ws.Range("E8").Value = 0.002183406113537118
$ws.Range("F8").Value = 0.03930131004366812
$ws.Range("J8").Value = 0.1135371179039301
$ws.Range("O8").Value = 0.01091703056768559
$ws.Range("Q8").Value = 0.2183406113537118
$ws.Range("R8").Value = 0.09606986899563319
$ws.Range("S8").Value = 0.4126637554585153

# Row 9
$ws.Range("B9").Value = 0.1036789297658863
$ws.Range("D9").Value = 0.02006688963210702
$ws.Range("F9").Value = 0.0568561872909699
$ws.Range("J9").Value = 0.1070234113712375
$ws.Range("O9").Value = 0.01003344481605351
$ws.Range("Q9").Value = 0.1538461538461539
$ws.Range("R9").Value = 0.1137123745819398
$ws.Range("S9").Value = 0.4347826086956522

# Row 10
$ws.Range("B10").Value = 0.1216320246343341
$ws.Range("D10").Value = 0.02155504234026174
$ws.Range("F10").Value = 0.07775211701308699
$ws.Range("J10").Value = 0.1054657428791378
$ws.Range("O10").Value = 0.0138568129330254
$ws.Range("Q10").Value = 0.1562740569668976
$ws.Range("R10").Value = 0.1062355658198614
$ws.Range("S10").Value = 0.3972286374133949

# Row 11
$ws.Range("G11").Value = 0.1302681992337165
$ws.Range("J11").Value = 0.08045977011494253
$ws.Range("K11").Value = 0.1954022988505747
$ws.Range("L11").Value = 0.5823754789272031
$ws.Range("S11").Value = 0.01149425287356322

# Row 12
$ws.Range("G12").Value = 0.8145695364238411
$ws.Range("J12").Value = 0.1721854304635762
$ws.Range("S12").Value = 0.01324503311258278

# Row 13
$ws.Range("G13").Value = 0.5757575757575758
$ws.Range("J13").Value = 0.3636363636363636
$ws.Range("S13").Value = 0.06060606060606061

# Row 15
$ws.Range("F15").Value = 0.01968503937007874
$ws.Range("H15").Value = 0.1653543307086614
$ws.Range("I15").Value = 0.09448818897637795
$ws.Range("J15").Value = 0.3307086614173229
$ws.Range("K15").Value = 0.03937007874015748
$ws.Range("M15").Value = 0.01968503937007874
$ws.Range("O15").Value = 0.09842519685039371
$ws.Range("S15").Value = 0.2322834645669291

# Row 16
$ws.Range("F16").Value = 0.04524886877828054
$ws.Range("H16").Value = 0.1312217194570136
$ws.Range("I16").Value = 0.1266968325791855
$ws.Range("J16").Value = 0.3755656108597285
$ws.Range("K16").Value = 0.09502262443438914
$ws.Range("N16").Value = 0.004524886877828055
$ws.Range("O16").Value = 0.07692307692307693
$ws.Range("S16").Value = 0.1447963800904978

# Row 17
$ws.Range("F17").Value = 0.02179176755447942
$ws.Range("H17").Value = 0.1525423728813559
$ws.Range("I17").Value = 0.1404358353510896
$ws.Range("J17").Value = 0.3922518159806295
$ws.Range("K17").Value = 0.07263922518159806
$ws.Range("M17").Value = 0.01937046004842615
$ws.Range("N17").Value = 0.002421307506053269
$ws.Range("O17").Value = 0.06295399515738499
$ws.Range("S17").Value = 0.1355932203389831

# Row 18
$ws.Range("F18").Value = 0.02510460251046025
$ws.Range("H18").Value = 0.1464435146443515
$ws.Range("I18").Value = 0.1380753138075314
$ws.Range("J18").Value = 0.3723849372384937
$ws.Range("K18").Value = 0.08786610878661087
$ws.Range("M18").Value = 0.008368200836820083
$ws.Range("O18").Value = 0.08368200836820083
$ws.Range("S18").Value = 0.1380753138075314

# Row 19
$ws.Range("F19").Value = 0.02110625909752547
$ws.Range("H19").Value = 0.2132459970887919
$ws.Range("I19").Value = 0.1157205240174673
$ws.Range("J19").Value = 0.3646288209606987
$ws.Range("K19").Value = 0.08879184861717612
$ws.Range("M19").Value = 0.01382823871906841
$ws.Range("N19").Value = 0.000727802037845706
$ws.Range("O19").Value = 0.07641921397379912
$ws.Range("S19").Value = 0.1055312954876274

Write-Host "Applied team-specific matrix updates"
